$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5796
$ws.Range("I116").Value = 3892
$ws.Range("J116").Value = 6661.4546
$ws.Range("K116").Value = 3892
$ws.Range("L116").Value = 6661.4546
$ws.Range("M116").Value = -450
$ws.Range("N116").Value = -13545.4546

$ws.Range("H127").Value = 924.94116
$ws.Range("I127").Value = 317.4
$ws.Range("J127").Value = 1792.8572
$ws.Range("K127").Value = 952.1999999999999
$ws.Range("L127").Value = 5378.571599999999
$ws.Range("M127").Value = 4007.8
$ws.Range("N127").Value = -15298.5716

$ws.Range("H137").Value = 286981.53
$ws.Range("I137").Value = 440765.7
$ws.Range("J137").Value = 3694.8948
$ws.Range("K137").Value = 1322297.1
$ws.Range("L137").Value = 11084.6844
$ws.Range("M137").Value = -1319747.1
$ws.Range("N137").Value = -16184.6844

$ws.Range("H138").Value = 3059.03
$ws.Range("I138").Value = 1642.3214
$ws.Range("J138").Value = 3609.9722
$ws.Range("K138").Value = 4926.9642
$ws.Range("L138").Value = 10829.9166
$ws.Range("M138").Value = 213.0357999999997
$ws.Range("N138").Value = -21109.9166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 26048
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 26048
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 26048
$ws.Range("N52").Value = -26684

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 29314
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 29314
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 29314
$ws.Range("N58").Value = -29902

$ws.Range("H86").Value = 1726.2903
$ws.Range("I86").Value = 1802.5883
$ws.Range("J86").Value = 1633.6428
$ws.Range("K86").Value = 1802.5883
$ws.Range("L86").Value = 1633.6428
$ws.Range("M86").Value = -679.5882999999999
$ws.Range("N86").Value = -3879.6428

$ws.Range("H89").Value = 1726.2903
$ws.Range("I89").Value = 1802.5883
$ws.Range("J89").Value = 1633.6428
$ws.Range("K89").Value = 9012.941499999999
$ws.Range("L89").Value = 8168.214
$ws.Range("M89").Value = -3396.941499999999
$ws.Range("N89").Value = -19400.214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.5823
$ws.Range("I31").Value = 1193.5869
$ws.Range("J31").Value = 4707.879
$ws.Range("K31").Value = 1193.5869
$ws.Range("L31").Value = 4707.879
$ws.Range("M31").Value = -898.5869
$ws.Range("N31").Value = -5297.879

$ws.Range("H34").Value = 2661.5823
$ws.Range("I34").Value = 1193.5869
$ws.Range("J34").Value = 4707.879
$ws.Range("K34").Value = 1193.5869
$ws.Range("L34").Value = 4707.879
$ws.Range("M34").Value = -991.5869
$ws.Range("N34").Value = -5111.879

$ws.Range("H99").Value = 1707.4546
$ws.Range("I99").Value = 1359.8
$ws.Range("J99").Value = 1997.1666
$ws.Range("K99").Value = 1359.8
$ws.Range("L99").Value = 1997.1666
$ws.Range("M99").Value = 138.2
$ws.Range("N99").Value = -4993.1666

$ws.Range("H105").Value = 1515.8975
$ws.Range("I105").Value = 966.2069
$ws.Range("J105").Value = 3110
$ws.Range("K105").Value = 966.2069
$ws.Range("L105").Value = 3110
$ws.Range("M105").Value = 780.7931
$ws.Range("N105").Value = -6604

$ws.Range("H122").Value = 2580
$ws.Range("I122").Value = 1433.3334
$ws.Range("J122").Value = 4300
$ws.Range("K122").Value = 4300.0002
$ws.Range("L122").Value = 12900
$ws.Range("M122").Value = -1850.0002
$ws.Range("N122").Value = -17800

$ws.Range("H126").Value = 1707.4546
$ws.Range("I126").Value = 1359.8
$ws.Range("J126").Value = 1997.1666
$ws.Range("K126").Value = 4079.4
$ws.Range("L126").Value = 5991.4998
$ws.Range("M126").Value = -1609.4
$ws.Range("N126").Value = -10931.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61
$ws.Range("I2").Value = 58.75
$ws.Range("J2").Value = 64
$ws.Range("K2").Value = 58.75
$ws.Range("L2").Value = 64
$ws.Range("M2").Value = 54.25
$ws.Range("N2").Value = -290

$ws.Range("H15").Value = 33500
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 33500
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 33500
$ws.Range("N15").Value = -34076

$ws.Range("H41").Value = 9200
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 9200
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 9200
$ws.Range("N41").Value = -9910

$ws.Range("H81").Value = 33500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 33500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 33500
$ws.Range("N81").Value = -35496

$ws.Range("H84").Value = 33500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 33500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 100500
$ws.Range("N84").Value = -110484

$ws.Range("H122").Value = 2063.5483
$ws.Range("I122").Value = 1983.7778
$ws.Range("J122").Value = 2602
$ws.Range("K122").Value = 5951.3334
$ws.Range("L122").Value = 7806
$ws.Range("M122").Value = -3501.3334
$ws.Range("N122").Value = -12706

$ws.Range("H126").Value = 12296.173
$ws.Range("I126").Value = 3913.5
$ws.Range("J126").Value = 20120
$ws.Range("K126").Value = 11740.5
$ws.Range("L126").Value = 60360
$ws.Range("M126").Value = -9270.5
$ws.Range("N126").Value = -65300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25003050
$ws.Range("I7").Value = 50002100
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 50002100
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -50001988
$ws.Range("N7").Value = -4224

$ws.Range("H40").Value = 3525.5715
$ws.Range("I40").Value = 3717.75
$ws.Range("J40").Value = 2372.5
$ws.Range("K40").Value = 3717.75
$ws.Range("L40").Value = 2372.5
$ws.Range("M40").Value = -3581.75
$ws.Range("N40").Value = -2644.5

$ws.Range("H55").Value = 537.7143
$ws.Range("I55").Value = 116.23077
$ws.Range("J55").Value = 903
$ws.Range("K55").Value = 116.23077
$ws.Range("L55").Value = 903
$ws.Range("M55").Value = 56.76922999999999
$ws.Range("N55").Value = -1249

$ws.Range("H122").Value = 4331.591
$ws.Range("I122").Value = 3916.6667
$ws.Range("J122").Value = 4829.5
$ws.Range("K122").Value = 11750.0001
$ws.Range("L122").Value = 14488.5
$ws.Range("M122").Value = -9300.000100000001
$ws.Range("N122").Value = -19388.5

$ws.Range("H126").Value = 25003050
$ws.Range("I126").Value = 50002100
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 150006300
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -150003830
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 4484
$ws.Range("I132").Value = 5071.3335
$ws.Range("J132").Value = 3309.3333
$ws.Range("K132").Value = 15214.0005
$ws.Range("L132").Value = 9927.999899999999
$ws.Range("M132").Value = -12684.0005
$ws.Range("N132").Value = -14987.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10192.667
$ws.Range("I122").Value = 13853.353
$ws.Range("J122").Value = 3969.5
$ws.Range("K122").Value = 41560.05899999999
$ws.Range("L122").Value = 11908.5
$ws.Range("M122").Value = -39110.05899999999
$ws.Range("N122").Value = -16808.5

$ws.Range("H127").Value = 33929.8
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 33929.8
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 33929.8
$ws.Range("N127").Value = -43849.8
